$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "VALOR MORA" total and "Cant. Trabajadores" now that the
# second worker (VICTOR RAUL MENDOZA MARTINEZ) is removed.
$ws.Range("E11").Value = 250753
$ws.Range("C13").Value = 1

# Rewrite the remaining worker's (LUIS CARLOS CARCAMO ROMAN) period rows in
# ascending period order, swapping the F-column (Valor Mora) values for the
# first/last period so each period keeps its original amount.
$ws.Range("E16").Value = "2102"
$ws.Range("F16").Value = 32707

$ws.Range("E17").Value = "2103"
$ws.Range("F17").Value = 36341

$ws.Range("E18").Value = "2104"
$ws.Range("F18").Value = 36341

$ws.Range("E19").Value = "2105"
$ws.Range("F19").Value = 36341

$ws.Range("E20").Value = "2106"
$ws.Range("F20").Value = 36341

$ws.Range("E21").Value = "2107"
$ws.Range("F21").Value = 36341

$ws.Range("E22").Value = "2108"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

# Give row 22 the closing (bottom-border) style that used to belong to the
# last data row (25), since it is now the last row of the table.
$ws.Range("B25:J25").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Delete the three rows that belonged to the removed worker
# (VICTOR RAUL MENDOZA MARTINEZ, 73072871) - old rows 23, 24, 25.
$ws.Range("A23:J25").Delete()

# The signature rows move up from 30/31 to 27/28 as part of the same
# deletion (handled automatically by the row delete above since everything
# below shifts up).

# Narrow column D now that the longer "VICTOR RAUL MENDOZA MARTINEZ" text is
# gone and "LUIS CARLOS CARCAMO ROMAN" is the widest remaining entry.
$ws.Columns("D").ColumnWidth = 29
